$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-added parameter values in column B
$ws.Range("B12").Value = 0.6147
$ws.Range("B13").Value = 0.037
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 0.2
$ws.Range("B17").Value = 10
$ws.Range("B18").Value = 9
$ws.Range("B19").Value = 2

# Update the view: scroll back to top-left A1 and move selection to B19
$ws.Range("A1").Select() | Out-Null
$ws.Range("B19").Select() | Out-Null
